# Atualização de bases das ligas, do dia: 09-04-2024 às 22:40
#
# The underlying data rows got reshuffled: for a handful of row pairs the
# entire record (every column except the running "id" in column A) was
# swapped with its neighbour, and one trio of rows got rotated. Column A
# (the sequential id) and column E (match date, identical within each
# swapped group) are left untouched; everything from column B ("id"/match
# id in the source feed) through column AC (PL_AhUnder) moves together so
# the row keeps referring to one real match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $rowA, $rowB) {
    $rngA = $ws.Range("B$rowA`:AC$rowA")
    $rngB = $ws.Range("B$rowB`:AC$rowB")
    $valA = $rngA.Value2
    $valB = $rngB.Value2
    $rngA.Value2 = $valB
    $rngB.Value2 = $valA
}

# Simple pairwise swaps (each pair exchanges all of B:AC).
$pairs = @(
    @(23, 24),
    @(32, 33),
    @(60, 61),
    @(74, 75),
    @(132, 133),
    @(134, 135),
    @(140, 141),
    @(142, 143),
    @(151, 152),
    @(201, 202),
    @(230, 231),
    @(243, 244),
    @(247, 248)
)

foreach ($pair in $pairs) {
    Swap-Rows $ws $pair[0] $pair[1]
}

# Three-way rotation: row 221 takes row 222's old content, row 222 takes
# row 223's old content, and row 223 takes row 221's old content.
$r221 = $ws.Range("B221:AC221")
$r222 = $ws.Range("B222:AC222")
$r223 = $ws.Range("B223:AC223")

$v221 = $r221.Value2
$v222 = $r222.Value2
$v223 = $r223.Value2

$r221.Value2 = $v222
$r222.Value2 = $v223
$r223.Value2 = $v221
